$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 261, shifting rows 261:393 down to 262:394
$ws.Rows.Item(261).Insert()

# The columns A,B,C,E,F,G,H,I,N,O,Q,R keep the same constant values as the
# rest of the dataset; copy them from the row that is now at 262.
$ws.Range("A261").Value2 = $ws.Range("A262").Value2
$ws.Range("B261").Value2 = $ws.Range("B262").Value2
$ws.Range("C261").Value2 = $ws.Range("C262").Value2
$ws.Range("E261").Value2 = $ws.Range("E262").Value2
$ws.Range("F261").Value2 = $ws.Range("F262").Value2
$ws.Range("G261").Value2 = $ws.Range("G262").Value2
$ws.Range("H261").Value2 = $ws.Range("H262").Value2
$ws.Range("I261").Value2 = $ws.Range("I262").Value2
$ws.Range("N261").Value2 = $ws.Range("N262").Value2
$ws.Range("O261").Value2 = $ws.Range("O262").Value2
$ws.Range("Q261").Value2 = $ws.Range("Q262").Value2
$ws.Range("R261").Value2 = $ws.Range("R262").Value2

# New data values for the inserted row
$ws.Range("D261").Value2 = 44813
$ws.Range("J261").Value2 = 500
$ws.Range("K261").Value2 = 10000
$ws.Range("L261").Value2 = 10000
$ws.Range("M261").Value2 = 10000
$ws.Range("P261").Value2 = 500
